# Update "想去人数" (column F) counts across the four sheets to match the
# refreshed scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 19
$ws.Cells.Item(6, 6).Value = 1066
$ws.Cells.Item(8, 6).Value = 151
$ws.Cells.Item(9, 6).Value = 536
$ws.Cells.Item(10, 6).Value = 41
$ws.Cells.Item(11, 6).Value = 432
$ws.Cells.Item(12, 6).Value = 167
$ws.Cells.Item(13, 6).Value = 1288
$ws.Cells.Item(14, 6).Value = 1173
$ws.Cells.Item(15, 6).Value = 1354
$ws.Cells.Item(16, 6).Value = 8
$ws.Cells.Item(17, 6).Value = 16
$ws.Cells.Item(18, 6).Value = 264
$ws.Cells.Item(19, 6).Value = 1510
$ws.Cells.Item(22, 6).Value = 291
$ws.Cells.Item(25, 6).Value = 1059
$ws.Cells.Item(26, 6).Value = 295
$ws.Cells.Item(27, 6).Value = 777
$ws.Cells.Item(29, 6).Value = 937
$ws.Cells.Item(30, 6).Value = 185400
$ws.Cells.Item(31, 6).Value = 909
$ws.Cells.Item(32, 6).Value = 544
$ws.Cells.Item(33, 6).Value = 1304
$ws.Cells.Item(35, 6).Value = 5
$ws.Cells.Item(36, 6).Value = 5
$ws.Cells.Item(37, 6).Value = 800
$ws.Cells.Item(38, 6).Value = 1513
$ws.Cells.Item(40, 6).Value = 7
$ws.Cells.Item(43, 6).Value = 749
$ws.Cells.Item(44, 6).Value = 96
$ws.Cells.Item(45, 6).Value = 23

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 103
$ws.Cells.Item(8, 6).Value = 79
$ws.Cells.Item(11, 6).Value = 1339
$ws.Cells.Item(13, 6).Value = 2421
$ws.Cells.Item(14, 6).Value = 1151
$ws.Cells.Item(16, 6).Value = 706
$ws.Cells.Item(17, 6).Value = 177
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(28, 6).Value = 31
$ws.Cells.Item(30, 6).Value = 221
$ws.Cells.Item(33, 6).Value = 57
$ws.Cells.Item(37, 6).Value = 153
$ws.Cells.Item(42, 6).Value = 104

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 854
$ws.Cells.Item(6, 6).Value = 2647
$ws.Cells.Item(7, 6).Value = 4400
$ws.Cells.Item(8, 6).Value = 108
$ws.Cells.Item(10, 6).Value = 480
$ws.Cells.Item(11, 6).Value = 517
$ws.Cells.Item(12, 6).Value = 364
$ws.Cells.Item(13, 6).Value = 21
$ws.Cells.Item(14, 6).Value = 456
$ws.Cells.Item(15, 6).Value = 148

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 854
$ws.Cells.Item(5, 6).Value = 4400
$ws.Cells.Item(6, 6).Value = 108
$ws.Cells.Item(7, 6).Value = 480
$ws.Cells.Item(9, 6).Value = 456
$ws.Cells.Item(10, 6).Value = 148
$ws.Cells.Item(12, 6).Value = 1066
$ws.Cells.Item(14, 6).Value = 151
$ws.Cells.Item(16, 6).Value = 1339
$ws.Cells.Item(17, 6).Value = 536
$ws.Cells.Item(18, 6).Value = 41
$ws.Cells.Item(19, 6).Value = 432
$ws.Cells.Item(20, 6).Value = 167
$ws.Cells.Item(21, 6).Value = 2421
$ws.Cells.Item(22, 6).Value = 1151
$ws.Cells.Item(23, 6).Value = 1288
$ws.Cells.Item(24, 6).Value = 1173
$ws.Cells.Item(25, 6).Value = 1355
$ws.Cells.Item(26, 6).Value = 706
$ws.Cells.Item(27, 6).Value = 177
$ws.Cells.Item(28, 6).Value = 264
$ws.Cells.Item(29, 6).Value = 55
$ws.Cells.Item(30, 6).Value = 1511
$ws.Cells.Item(32, 6).Value = 291
$ws.Cells.Item(34, 6).Value = 1059
$ws.Cells.Item(35, 6).Value = 777
$ws.Cells.Item(37, 6).Value = 937
$ws.Cells.Item(39, 6).Value = 910
$ws.Cells.Item(40, 6).Value = 544
$ws.Cells.Item(42, 6).Value = 800
$ws.Cells.Item(43, 6).Value = 1513
$ws.Cells.Item(45, 6).Value = 153
$ws.Cells.Item(48, 6).Value = 749
$ws.Cells.Item(49, 6).Value = 96
